$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
# Row 14 (UPDATE) now marked Done
$ws.Range("C14").Value = 1

# Row 15: CREATE TABLE -> CREATE TABLE/VIEW, category DB management utilities -> ETL procedures, marked Done
$ws.Range("A15").Value = "CREATE TABLE/VIEW"
$ws.Range("B15").Value = "ETL procedures"
$ws.Range("C15").Value = 1

# Row 16 (INSERT) now marked Done
$ws.Range("C16").Value = 1

# --- Header for new column C ---
$ws.Range("C1").Value = "Done"

# --- AutoFilter over the table range (applied before appending new rows so
#     the filter range stays pinned to A1:C16, matching the target) ---
$ws.Range("A1:C16").AutoFilter(1)

# Register the hidden _FilterDatabase defined name scoped to this sheet
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Hoja1!`$A`$1:`$C`$16")
$fd.Visible = $false

# --- Append new rows ---
# Row 17: DELETE
$ws.Range("A17").Value = "DELETE"
$ws.Range("B17").Value = "ETL procedures"
$ws.Range("C17").Value = 1

# Row 18: DATA TYPES
$ws.Range("A18").Value = "DATA TYPES"
$ws.Range("B18").Value = "Data types (ORACLE SQL)"

# Row 19: INTEGRITY CONSTRAINTS
$ws.Range("A19").Value = "INTEGRITY CONSTRAINTS"
$ws.Range("B19").Value = "DB management utilities"

# --- Sheet view: scroll so row 7 is at the top, select C17 ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
